# Generate Report for Handback
# This script updates the localization-status workbook so that files that
# were handed off now show as handed back (in sync with en-US). For each
# handed-off row on the per-language sheets (zh-cn, de-de) we now also
# record the Latest Target File / Latest Handback File (columns E/F) and
# stamp the Latest Handback DateTime (column G) with the handback time.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: the Status column for both languages reflects the
#     same underlying shared string, so update it there too. ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

$zh.Range("E2").Value = "8681868a-e100-4edc-8062-ccaeb78afc2b.md"
$zh.Range("F2").Value = "8681868a-e100-4edc-8062-ccaeb78afc2b.b9158bee296fb9e0e67d71e0cead111bc7bbd29c.zh-cn.xlf"
$zh.Range("G2").Value = "2016-03-02 10:37:15"

$zh.Range("E3").Value = "f2a91ce9-c14c-402b-8b0d-615309079abd.md"
$zh.Range("F3").Value = "f2a91ce9-c14c-402b-8b0d-615309079abd.fd76faf5b7d3691262ae8781a664b1f1f20afd28.zh-cn.xlf"
$zh.Range("G3").Value = "2016-03-02 10:37:15"

$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/30c0d61a1886b98cdce0066d6993dce32509b37d/e2e/8681868a-e100-4edc-8062-ccaeb78afc2b.md", [Type]::Missing, [Type]::Missing, "8681868a-e100-4edc-8062-ccaeb78afc2b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c38f9aff167d8848d31dae8f986fd4a632dcd2d8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8681868a-e100-4edc-8062-ccaeb78afc2b.b9158bee296fb9e0e67d71e0cead111bc7bbd29c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "8681868a-e100-4edc-8062-ccaeb78afc2b.b9158bee296fb9e0e67d71e0cead111bc7bbd29c.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/30c0d61a1886b98cdce0066d6993dce32509b37d/e2e/f2a91ce9-c14c-402b-8b0d-615309079abd.md", [Type]::Missing, [Type]::Missing, "f2a91ce9-c14c-402b-8b0d-615309079abd.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c38f9aff167d8848d31dae8f986fd4a632dcd2d8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f2a91ce9-c14c-402b-8b0d-615309079abd.fd76faf5b7d3691262ae8781a664b1f1f20afd28.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "f2a91ce9-c14c-402b-8b0d-615309079abd.fd76faf5b7d3691262ae8781a664b1f1f20afd28.zh-cn.xlf") | Out-Null

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

$de.Range("E2").Value = "8681868a-e100-4edc-8062-ccaeb78afc2b.md"
$de.Range("F2").Value = "8681868a-e100-4edc-8062-ccaeb78afc2b.b9158bee296fb9e0e67d71e0cead111bc7bbd29c.de-de.xlf"
$de.Range("G2").Value = "2016-03-02 10:37:35"

$de.Range("E3").Value = "f2a91ce9-c14c-402b-8b0d-615309079abd.md"
$de.Range("F3").Value = "f2a91ce9-c14c-402b-8b0d-615309079abd.fd76faf5b7d3691262ae8781a664b1f1f20afd28.de-de.xlf"
$de.Range("G3").Value = "2016-03-02 10:37:35"

$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/30c0d61a1886b98cdce0066d6993dce32509b37d/e2e/8681868a-e100-4edc-8062-ccaeb78afc2b.md", [Type]::Missing, [Type]::Missing, "8681868a-e100-4edc-8062-ccaeb78afc2b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da1b9b859336df2f6832e654de89d29afadebe72/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8681868a-e100-4edc-8062-ccaeb78afc2b.b9158bee296fb9e0e67d71e0cead111bc7bbd29c.de-de.xlf", [Type]::Missing, [Type]::Missing, "8681868a-e100-4edc-8062-ccaeb78afc2b.b9158bee296fb9e0e67d71e0cead111bc7bbd29c.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/30c0d61a1886b98cdce0066d6993dce32509b37d/e2e/f2a91ce9-c14c-402b-8b0d-615309079abd.md", [Type]::Missing, [Type]::Missing, "f2a91ce9-c14c-402b-8b0d-615309079abd.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da1b9b859336df2f6832e654de89d29afadebe72/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f2a91ce9-c14c-402b-8b0d-615309079abd.fd76faf5b7d3691262ae8781a664b1f1f20afd28.de-de.xlf", [Type]::Missing, [Type]::Missing, "f2a91ce9-c14c-402b-8b0d-615309079abd.fd76faf5b7d3691262ae8781a664b1f1f20afd28.de-de.xlf") | Out-Null
